# Scheduled-runner refresh of market-price derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2683.1133
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2683.1133
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 8049.3399
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -8385.339899999999

$ws.Range("H62").Value = 5080.625
$ws.Range("I62").Value = 6132.222
$ws.Range("J62").Value = 3728.5715
$ws.Range("K62").Value = 6132.222
$ws.Range("L62").Value = 3728.5715
$ws.Range("M62").Value = -5508.222
$ws.Range("N62").Value = -4976.5715

$ws.Range("H65").Value = 5080.625
$ws.Range("I65").Value = 6132.222
$ws.Range("J65").Value = 3728.5715
$ws.Range("K65").Value = 30661.11
$ws.Range("L65").Value = 18642.8575
$ws.Range("M65").Value = -27541.11
$ws.Range("N65").Value = -24882.8575

$ws.Range("H133").Value = 58585.43
$ws.Range("J133").Value = 58585.43
$ws.Range("L133").Value = 58585.43
$ws.Range("N133").Value = -68705.42999999999

$ws.Range("H137").Value = 3185.2834
$ws.Range("I137").Value = 1001.1539
$ws.Range("J137").Value = 7241.524
$ws.Range("K137").Value = 3003.4617
$ws.Range("L137").Value = 21724.572
$ws.Range("M137").Value = -453.4616999999998
$ws.Range("N137").Value = -26824.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9813.984
$ws.Range("I32").Value = 8951.232
$ws.Range("J32").Value = 16716
$ws.Range("K32").Value = 8951.232
$ws.Range("L32").Value = 16716
$ws.Range("M32").Value = -8664.232
$ws.Range("N32").Value = -17290

$ws.Range("H45").Value = 1406.68
$ws.Range("I45").Value = 1234.8422
$ws.Range("J45").Value = 1950.8334
$ws.Range("K45").Value = 1234.8422
$ws.Range("L45").Value = 1950.8334
$ws.Range("M45").Value = -857.8422
$ws.Range("N45").Value = -2704.8334

$ws.Range("H64").Value = 29908.908
$ws.Range("J64").Value = 29908.908
$ws.Range("L64").Value = 29908.908
$ws.Range("N64").Value = -30404.908

$ws.Range("H67").Value = 29908.908
$ws.Range("J67").Value = 29908.908
$ws.Range("L67").Value = 29908.908
$ws.Range("N67").Value = -31624.908

$ws.Range("H95").Value = 40103.5
$ws.Range("J95").Value = 40103.5
$ws.Range("L95").Value = 40103.5
$ws.Range("N95").Value = -45595.5

$ws.Range("H105").Value = 49366
$ws.Range("J105").Value = 49366
$ws.Range("L105").Value = 49366
$ws.Range("N105").Value = -56354

$ws.Range("H123").Value = 41777.5
$ws.Range("J123").Value = 41777.5
$ws.Range("L123").Value = 41777.5
$ws.Range("N123").Value = -51577.5

$ws.Range("H138").Value = 53300
$ws.Range("J138").Value = 53300
$ws.Range("L138").Value = 53300
$ws.Range("N138").Value = -63580

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 58999.332
$ws.Range("J57").Value = 58999.332
$ws.Range("L57").Value = 58999.332
$ws.Range("N57").Value = -60439.332

$ws.Range("H62").Value = 30695.25
$ws.Range("J62").Value = 30695.25
$ws.Range("L62").Value = 30695.25
$ws.Range("N62").Value = -32067.25

$ws.Range("H65").Value = 30695.25
$ws.Range("J65").Value = 30695.25
$ws.Range("L65").Value = 92085.75
$ws.Range("N65").Value = -98949.75

$ws.Range("H122").Value = 40549.4
$ws.Range("J122").Value = 40549.4
$ws.Range("L122").Value = 40549.4
$ws.Range("N122").Value = -50349.4

$ws.Range("H132").Value = 24898.572
$ws.Range("J132").Value = 24898.572
$ws.Range("L132").Value = 24898.572
$ws.Range("N132").Value = -35018.572

$ws.Range("H133").Value = 40666.5
$ws.Range("J133").Value = 40666.5
$ws.Range("L133").Value = 40666.5
$ws.Range("N133").Value = -50786.5

$ws.Range("H136").Value = 58999.332
$ws.Range("J136").Value = 58999.332
$ws.Range("L136").Value = 58999.332
$ws.Range("N136").Value = -69199.33199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 48657
$ws.Range("J43").Value = 48657
$ws.Range("L43").Value = 48657
$ws.Range("N43").Value = -49025

$ws.Range("H52").Value = 76500
$ws.Range("J52").Value = 76500
$ws.Range("L52").Value = 76500
$ws.Range("N52").Value = -77088

$ws.Range("H96").Value = 71864.8
$ws.Range("J96").Value = 71864.8
$ws.Range("L96").Value = 71864.8
$ws.Range("N96").Value = -77356.8

$ws.Range("H101").Value = 48657
$ws.Range("J101").Value = 48657
$ws.Range("L101").Value = 48657
$ws.Range("N101").Value = -55147

$ws.Range("H105").Value = 2548.524
$ws.Range("I105").Value = 2514.0588
$ws.Range("K105").Value = 2514.0588
$ws.Range("M105").Value = -767.0587999999998

$ws.Range("H134").Value = 42678.44
$ws.Range("I134").Value = 1340.3448
$ws.Range("K134").Value = 4021.0344
$ws.Range("M134").Value = -1486.0344

$ws.Range("H137").Value = 45799.93
$ws.Range("J137").Value = 45799.93
$ws.Range("L137").Value = 45799.93
$ws.Range("N137").Value = -55999.93

$ws.Range("H139").Value = 62899.8
$ws.Range("J139").Value = 68124.75
$ws.Range("L139").Value = 68124.75
$ws.Range("N139").Value = -78404.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1359.8182
$ws.Range("I4").Value = 92.833336
$ws.Range("J4").Value = 2880.2
$ws.Range("K4").Value = 278.500008
$ws.Range("L4").Value = 8640.599999999999
$ws.Range("M4").Value = -166.500008
$ws.Range("N4").Value = -8864.599999999999

$ws.Range("H113").Value = 2719.3125
$ws.Range("I113").Value = 3728.9355
$ws.Range("J113").Value = 878.2353000000001
$ws.Range("K113").Value = 11186.8065
$ws.Range("L113").Value = 2634.7059
$ws.Range("M113").Value = -9016.806500000001
$ws.Range("N113").Value = -6974.7059

$ws.Range("H131").Value = 3871.0789
$ws.Range("J131").Value = 1524.0358
$ws.Range("L131").Value = 4572.107400000001
$ws.Range("N131").Value = -14652.1074

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H104").Value = 41210.168
$ws.Range("J104").Value = 41210.168
$ws.Range("L104").Value = 41210.168
$ws.Range("N104").Value = -48198.168

$ws.Range("H113").Value = 1300.9286
$ws.Range("I113").Value = 1371.4286
$ws.Range("J113").Value = 1230.4286
$ws.Range("K113").Value = 1371.4286
$ws.Range("L113").Value = 1230.4286
$ws.Range("M113").Value = 798.5714
$ws.Range("N113").Value = -5570.4286

$ws.Range("H135").Value = 43140
$ws.Range("J135").Value = 43140
$ws.Range("L135").Value = 43140
$ws.Range("N135").Value = -53280

$ws.Range("H138").Value = 53000
$ws.Range("J138").Value = 53000
$ws.Range("L138").Value = 53000
$ws.Range("N138").Value = -63280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 39671.8
$ws.Range("J110").Value = 39671.8
$ws.Range("L110").Value = 39671.8
$ws.Range("N110").Value = -47851.8

$ws.Range("H134").Value = 47968.332
$ws.Range("J134").Value = 47968.332
$ws.Range("L134").Value = 47968.332
$ws.Range("N134").Value = -58108.332

$ws.Range("H137").Value = 40800
$ws.Range("J137").Value = 40800
$ws.Range("L137").Value = 40800
$ws.Range("N137").Value = -51000

$ws.Range("H139").Value = 50899.75
$ws.Range("J139").Value = 50899.75
$ws.Range("L139").Value = 50899.75
$ws.Range("N139").Value = -61179.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 55138.75
$ws.Range("J46").Value = 55138.75
$ws.Range("L46").Value = 55138.75
$ws.Range("N46").Value = -55600.75

$ws.Range("H105").Value = 47307.5
$ws.Range("J105").Value = 47307.5
$ws.Range("L105").Value = 47307.5
$ws.Range("N105").Value = -54295.5

$ws.Range("H126").Value = 1961504.2
$ws.Range("I126").Value = 2942036.2
$ws.Range("J126").Value = 440
$ws.Range("K126").Value = 8826108.600000001
$ws.Range("L126").Value = 1320
$ws.Range("M126").Value = -8823638.600000001
$ws.Range("N126").Value = -6260

$ws.Range("H131").Value = 49058.8
$ws.Range("J131").Value = 49058.8
$ws.Range("L131").Value = 49058.8
$ws.Range("N131").Value = -59138.8

$ws.Range("H134").Value = 55138.75
$ws.Range("J134").Value = 55138.75
$ws.Range("L134").Value = 165416.25
$ws.Range("N134").Value = -170486.25

$ws.Range("H139").Value = 55699.8
$ws.Range("J139").Value = 55699.8
$ws.Range("L139").Value = 55699.8
$ws.Range("N139").Value = -65979.8
